$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2 through 54
# from Excel serial date 45174 to 45175 (one day later).
$ws.Range("C2:C54").Value = 45175
